# Refresh the per-coin Price (D) and Volume(1h) (E) columns with the latest
# crypto snapshot values. All cells in these two columns are stored as plain
# text in the workbook (prices use '.' as a thousands separator in several
# rows, e.g. "58.242.72", and volumes are padded percentage strings like
# "  +0.37%  "), so every assignment is prefixed with a literal leading
# apostrophe -- Excel's standard "force text" quote-prefix -- to stop
# clean-looking numbers ("546.92", "11.00", ...) from being silently
# reinterpreted as numeric values. The apostrophe itself is not stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''58.242.72'
$ws.Range('E2').Value = '''  +0.37%  '
$ws.Range('D3').Value = '''2.369.00'
$ws.Range('E3').Value = '''  +0.94%  '
$ws.Range('E4').Value = '''  +0.12%  '
$ws.Range('D5').Value = '''546.92'
$ws.Range('E5').Value = '''  +0.43%  '
$ws.Range('D6').Value = '''133.72'
$ws.Range('E6').Value = '''  -0.75%  '
$ws.Range('E7').Value = '''  +0.09%  '
$ws.Range('D8').Value = '''0.566'
$ws.Range('E8').Value = '''  +5.14%  '
$ws.Range('E9').Value = '''  +3.98%  '
$ws.Range('D10').Value = '''5.56'
$ws.Range('E10').Value = '''  +2.68%  '
$ws.Range('E11').Value = '''  -2.00%  '
$ws.Range('D12').Value = '''0.355'
$ws.Range('E12').Value = '''  -1.12%  '
$ws.Range('D13').Value = '''24.21'
$ws.Range('E13').Value = '''  +2.79%  '
$ws.Range('D14').Value = '''2.792.24'
$ws.Range('E14').Value = '''  +1.06%  '
$ws.Range('D15').Value = '''58.234.86'
$ws.Range('E15').Value = '''  +0.43%  '
$ws.Range('E16').Value = '''  +2.29%  '
$ws.Range('D17').Value = '''2.364.45'
$ws.Range('E17').Value = '''  +0.63%  '
$ws.Range('D18').Value = '''11.00'
$ws.Range('E18').Value = '''  +3.73%  '
$ws.Range('D19').Value = '''4.33'
$ws.Range('E19').Value = '''  +2.76%  '
$ws.Range('D20').Value = '''331.29'
$ws.Range('E20').Value = '''  -0.86%  '
$ws.Range('D21').Value = '''6.89'
$ws.Range('E21').Value = '''  +2.75%  '
$ws.Range('D22').Value = '''0.999'
$ws.Range('E22').Value = '''  +0.05%  '
$ws.Range('D23').Value = '''63.51'
$ws.Range('E23').Value = '''  +3.02%  '
$ws.Range('E24').Value = '''  -0.80%  '
$ws.Range('D25').Value = '''0.998'
$ws.Range('E25').Value = '''  -0.03%  '
$ws.Range('D26').Value = '''8.26'
$ws.Range('E26').Value = '''  -2.64%  '
$ws.Range('D27').Value = '''1.33'
$ws.Range('E27').Value = '''  -5.74%  '
$ws.Range('E28').Value = '''  +0.22%  '
$ws.Range('D29').Value = '''170.44'
$ws.Range('E29').Value = '''  +0.11%  '
$ws.Range('E30').Value = '''  +1.80%  '
$ws.Range('E31').Value = '''  +0.47%  '
$ws.Range('D32').Value = '''18.46'
$ws.Range('E32').Value = '''  +0.02%  '
$ws.Range('E33').Value = '''  -0.03%  '
$ws.Range('D34').Value = '''0.997'
$ws.Range('E34').Value = '''  -3.94%  '
$ws.Range('E35').Value = '''  +0.08%  '
$ws.Range('D36').Value = '''4.18'
$ws.Range('E36').Value = '''  +0.20%  '
$ws.Range('E37').Value = '''  -1.55%  '
$ws.Range('E38').Value = '''  -2.00%  '
$ws.Range('D39').Value = '''0.412'
$ws.Range('E39').Value = '''  +8.69%  '
$ws.Range('D40').Value = '''142.83'
$ws.Range('E40').Value = '''  -4.14%  '
$ws.Range('D41').Value = '''3.69'
$ws.Range('E41').Value = '''  +2.32%  '
$ws.Range('D42').Value = '''288.16'
$ws.Range('E42').Value = '''  +0.78%  '
$ws.Range('D43').Value = '''0.0951'
$ws.Range('E43').Value = '''  +2.81%  '
$ws.Range('D44').Value = '''0.0519'
$ws.Range('E44').Value = '''  +2.81%  '
$ws.Range('D45').Value = '''18.96'
$ws.Range('E45').Value = '''  -1.38%  '
$ws.Range('D46').Value = '''0.566'
$ws.Range('E46').Value = '''  +0.68%  '
$ws.Range('E47').Value = '''  +2.52%  '
$ws.Range('D48').Value = '''0.391'
$ws.Range('E48').Value = '''  +2.37%  '
$ws.Range('D49').Value = '''11.08'
$ws.Range('E49').Value = '''  +0.24%  '
$ws.Range('E50').Value = '''  +0.78%  '
$ws.Range('D51').Value = '''1.54'
$ws.Range('E51').Value = '''  +0.38%  '
